$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-31 Wednesday" "2026-01-01 Thursday"

Replace-Text "59×17=" "32×27="
Replace-Text "26×31=" "48×12="
Replace-Text "80×89=" "90×17="
Replace-Text "97×17=" "93×98="
Replace-Text "52×68=" "97×48="
Replace-Text "29×49=" "24×44="
Replace-Text "57×30=" "18×19="
Replace-Text "91×24=" "74×96="
Replace-Text "57×33=" "92×66="
Replace-Text "84×15=" "70×66="
Replace-Text "86×48=" "95×73="
Replace-Text "74×91=" "46×22="
Replace-Text "73×28=" "43×30="
Replace-Text "94×69=" "64×48="
Replace-Text "83×44=" "34×20="
Replace-Text "97×71=" "64×38="
Replace-Text "57×41=" "96×68="
Replace-Text "44×72=" "15×92="
Replace-Text "84×69=" "54×86="
Replace-Text "87×55=" "72×66="
Replace-Text "93×80=" "96×66="
Replace-Text "67×54=" "66×92="
Replace-Text "44×37=" "83×94="
Replace-Text "38×17=" "22×84="
Replace-Text "26×19=" "37×55="
